# Auto-generated script applying Moogle_Profits market-data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1719.2545
$ws.Range("J17").Value = 1719.2545
$ws.Range("L17").Value = 5157.7635
$ws.Range("N17").Value = -5493.7635

$ws.Range("H62").Value = 3115.4736
$ws.Range("I62").Value = 2791.4614
$ws.Range("J62").Value = 3817.5
$ws.Range("K62").Value = 2791.4614
$ws.Range("L62").Value = 3817.5
$ws.Range("M62").Value = -2167.4614
$ws.Range("N62").Value = -5065.5

$ws.Range("H65").Value = 3115.4736
$ws.Range("I65").Value = 2791.4614
$ws.Range("J65").Value = 3817.5
$ws.Range("K65").Value = 13957.307
$ws.Range("L65").Value = 19087.5
$ws.Range("M65").Value = -10837.307
$ws.Range("N65").Value = -25327.5

$ws.Range("H70").Value = 2794.8
$ws.Range("I70").Value = 2714
$ws.Range("J70").Value = 2983.3333
$ws.Range("K70").Value = 8142
$ws.Range("L70").Value = 8949.999899999999
$ws.Range("M70").Value = -7872
$ws.Range("N70").Value = -9489.999899999999

$ws.Range("H73").Value = 2794.8
$ws.Range("I73").Value = 2714
$ws.Range("J73").Value = 2983.3333
$ws.Range("K73").Value = 8142
$ws.Range("L73").Value = 8949.999899999999
$ws.Range("M73").Value = -7206
$ws.Range("N73").Value = -10821.9999

$ws.Range("H107").Value = 908.0909
$ws.Range("I107").Value = 891.375
$ws.Range("J107").Value = 952.6667
$ws.Range("K107").Value = 891.375
$ws.Range("L107").Value = 952.6667
$ws.Range("M107").Value = 1028.625
$ws.Range("N107").Value = -4792.6667

$ws.Range("H132").Value = 2324.0312
$ws.Range("I132").Value = 2127.5
$ws.Range("K132").Value = 6382.5
$ws.Range("M132").Value = -3852.5

$ws.Range("H135").Value = 3178.75
$ws.Range("J135").Value = 5300
$ws.Range("L135").Value = 47700
$ws.Range("N135").Value = -52770

$ws.Range("H137").Value = 2878.182
$ws.Range("I137").Value = 1597.1
$ws.Range("K137").Value = 4791.299999999999
$ws.Range("M137").Value = -2241.299999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6161.5156
$ws.Range("I32").Value = 2765.6667
$ws.Range("K32").Value = 2765.6667
$ws.Range("M32").Value = -2478.6667

$ws.Range("H45").Value = 2031.2069
$ws.Range("I45").Value = 1788.4
$ws.Range("K45").Value = 1788.4
$ws.Range("M45").Value = -1411.4

$ws.Range("H97").Value = 1459.9412
$ws.Range("I97").Value = 1459.9412
$ws.Range("K97").Value = 1459.9412
$ws.Range("M97").Value = -963.9412

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2937.3333
$ws.Range("I134").Value = 2431.6365
$ws.Range("J134").Value = 8500
$ws.Range("K134").Value = 7294.9095
$ws.Range("L134").Value = 25500
$ws.Range("M134").Value = -4759.9095
$ws.Range("N134").Value = -30570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6682
$ws.Range("I31").Value = 2739.15
$ws.Range("J31").Value = 19824.834
$ws.Range("K31").Value = 2739.15
$ws.Range("L31").Value = 19824.834
$ws.Range("M31").Value = -2444.15
$ws.Range("N31").Value = -20414.834

$ws.Range("H34").Value = 6682
$ws.Range("I34").Value = 2739.15
$ws.Range("J34").Value = 19824.834
$ws.Range("K34").Value = 2739.15
$ws.Range("L34").Value = 19824.834
$ws.Range("M34").Value = -2537.15
$ws.Range("N34").Value = -20228.834

$ws.Range("H99").Value = 2977.1667
$ws.Range("I99").Value = 2926.9656
$ws.Range("J99").Value = 3185.1428
$ws.Range("K99").Value = 2926.9656
$ws.Range("L99").Value = 3185.1428
$ws.Range("M99").Value = -1428.9656
$ws.Range("N99").Value = -6181.1428

$ws.Range("H126").Value = 2977.1667
$ws.Range("I126").Value = 2926.9656
$ws.Range("J126").Value = 3185.1428
$ws.Range("K126").Value = 8780.8968
$ws.Range("L126").Value = 9555.428400000001
$ws.Range("M126").Value = -6310.8968
$ws.Range("N126").Value = -14495.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 33999
$ws.Range("J26").Value = 33999
$ws.Range("L26").Value = 33999
$ws.Range("N26").Value = -34559

$ws.Range("H50").Value = 33999
$ws.Range("J50").Value = 33999
$ws.Range("L50").Value = 33999
$ws.Range("N50").Value = -34995

$ws.Range("H80").Value = 8781.172
$ws.Range("I80").Value = 6926.4
$ws.Range("J80").Value = 13418.1
$ws.Range("K80").Value = 6926.4
$ws.Range("L80").Value = 13418.1
$ws.Range("M80").Value = -5928.4
$ws.Range("N80").Value = -15414.1

$ws.Range("H83").Value = 8781.172
$ws.Range("I83").Value = 6926.4
$ws.Range("J83").Value = 13418.1
$ws.Range("K83").Value = 34632
$ws.Range("L83").Value = 67090.5
$ws.Range("M83").Value = -29640
$ws.Range("N83").Value = -77074.5

$ws.Range("H132").Value = 3828.7083
$ws.Range("I132").Value = 2822.7856
$ws.Range("J132").Value = 5237
$ws.Range("K132").Value = 8468.356800000001
$ws.Range("L132").Value = 15711
$ws.Range("M132").Value = -5938.356800000001
$ws.Range("N132").Value = -20771

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2600.8096
$ws.Range("I16").Value = 2255.2666
$ws.Range("J16").Value = 3464.6667
$ws.Range("K16").Value = 2255.2666
$ws.Range("L16").Value = 3464.6667
$ws.Range("M16").Value = -2085.2666
$ws.Range("N16").Value = -3804.6667

$ws.Range("H56").Value = 31248.5
$ws.Range("I56").Value = 12498
$ws.Range("J56").Value = 49999
$ws.Range("K56").Value = 12498
$ws.Range("L56").Value = 49999
$ws.Range("M56").Value = -11807
$ws.Range("N56").Value = -51381

$ws.Range("H68").Value = 8250
$ws.Range("I68").Value = 9000
$ws.Range("K68").Value = 9000
$ws.Range("M68").Value = -8251

$ws.Range("H71").Value = 8250
$ws.Range("I71").Value = 9000
$ws.Range("K71").Value = 45000
$ws.Range("M71").Value = -41256

$ws.Range("H82").Value = 1038.6666
$ws.Range("I82").Value = 968.6
$ws.Range("K82").Value = 968.6
$ws.Range("M82").Value = -607.6

$ws.Range("H85").Value = 1038.6666
$ws.Range("I85").Value = 968.6
$ws.Range("K85").Value = 968.6
$ws.Range("M85").Value = 279.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 21176.23
$ws.Range("I96").Value = 2155.25
$ws.Range("J96").Value = 29630
$ws.Range("K96").Value = 2155.25
$ws.Range("L96").Value = 29630
$ws.Range("M96").Value = -782.25
$ws.Range("N96").Value = -32376

$ws.Range("H122").Value = 2238.5925
$ws.Range("I122").Value = 2128.2354
$ws.Range("J122").Value = 2426.2
$ws.Range("K122").Value = 6384.706200000001
$ws.Range("L122").Value = 7278.599999999999
$ws.Range("M122").Value = -3934.706200000001
$ws.Range("N122").Value = -12178.6

$ws.Range("H126").Value = 1505.0385
$ws.Range("I126").Value = 1369.5217
$ws.Range("K126").Value = 4108.5651
$ws.Range("M126").Value = -1638.5651

$ws.Range("H132").Value = 2981.7144
$ws.Range("I132").Value = 2987.8333
$ws.Range("J132").Value = 2945
$ws.Range("K132").Value = 8963.499899999999
$ws.Range("L132").Value = 8835
$ws.Range("M132").Value = -6433.499899999999
$ws.Range("N132").Value = -13895
